# "finestra incidenza 7gg centrata su ultimo g"
#
# Recenters the 7-day rolling-sum window (column C = "somma mobile 7gg.",
# column D = "somma mobile 7gg. per 100mila abitanti") so that the window
# ending on the CURRENT row (days r-6 .. r) is used, instead of the window
# that previously started on the current row (days r .. r+6).
#
# Consequences of the shift:
#   - The first 3 rows that used to have a computed value (rows 5-7) no
#     longer have 6 days of history before them, so they become blank.
#   - The last 3 rows (182-184) now have 6 full days of history before them,
#     so they get a computed value for the first time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$popolazione = 3199

# First row containing daily data (row 1 is the header row).
$firstDataRow = 2

# Last row containing daily data - found dynamically via End(xlUp) from the
# bottom of column A (xlUp = -4162).
$lastDataRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {

    $windowStart = $r - 6

    if ($windowStart -lt $firstDataRow) {
        # Not enough history yet (fewer than 7 days available) - blank out,
        # but only touch cells that actually still hold a computed value
        # (leave already-blank cells untouched).
        if ($ws.Cells.Item($r, 3).Value2 -ne "") {
            $ws.Cells.Item($r, 3).ClearContents()
        }
        if ($ws.Cells.Item($r, 4).Value2 -ne "") {
            $ws.Cells.Item($r, 4).ClearContents()
        }
        continue
    }

    $nuoviPos = 0
    for ($j = $windowStart; $j -le $r; $j++) {
        $nuoviPos = $nuoviPos + $ws.Cells.Item($j, 2).Value2
    }

    $ws.Cells.Item($r, 3).Value = $nuoviPos
    $ws.Cells.Item($r, 4).Value = $nuoviPos * 100000 / $popolazione
}
